# CropScapeConversions.xlsx update
# - Add a "Developed/Mixed" description for code 125 (row 127)
# - Append three new code/description rows at the bottom of the table:
#     code 255 -> Corn/Soybean Rotation
#     code 256 -> Other corn rotation
#     code 257 -> Other inconsistent
# - Highlight the newly-described cells with a yellow fill
# - Move the selection to reflect where the user ended up (B260)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new rows after the existing last row (row 256 / code 254).
$ws.Range("A257").Value = 255
$ws.Range("B257").Value = "Corn/Soybean Rotation"
$ws.Range("B257").Interior.Color = 65535

# Fill in the previously-blank description for code 125 (row 127, since row 1 is the header).
$ws.Range("B127").Value = "Developed/Mixed"
$ws.Range("B127").Interior.Color = 65535

$ws.Range("A258").Value = 256
$ws.Range("B258").Value = "Other corn rotation"
$ws.Range("B258").Interior.Color = 65535

$ws.Range("A259").Value = 257
$ws.Range("B259").Value = "Other inconsistent"
$ws.Range("B259").Interior.Color = 65535

# Reflect the final cursor/scroll position used while editing.
$ws.Range("B260").Select() | Out-Null
